# Update cryptos list with latest price/volume data (scraped Mon Nov 11 10:00:04 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "81.763.74"
$ws.Range("E2").Value = "  +2.97%  "
$ws.Range("D3").Value = "3.178.13"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'210.90"
$ws.Range("E5").Value = "  +2.80%  "
$ws.Range("D6").Value = "'622.16"
$ws.Range("E6").Value = "  -1.87%  "
$ws.Range("D7").Value = "'0.283"
$ws.Range("E7").Value = "  +20.71%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.583"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "3.169.50"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").Value = "'0.587"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").Value = "'0.0000253"
$ws.Range("E12").Value = "  +11.49%  "
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").Value = "'5.31"
$ws.Range("E14").Value = "  -3.88%  "
$ws.Range("D15").Value = "3.751.80"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "'31.48"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "81.427.44"
$ws.Range("E17").Value = "  +2.67%  "
$ws.Range("D18").Value = "3.167.49"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").Value = "'3.18"
$ws.Range("E19").Value = "  +3.75%  "
$ws.Range("D20").Value = "'13.97"
$ws.Range("E20").Value = "  -3.52%  "
$ws.Range("D21").Value = "'433.66"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("D22").Value = "'8.88"
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("D23").Value = "'5.08"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").Value = "'7.25"
$ws.Range("E24").Value = "  +5.61%  "
$ws.Range("D25").Value = "'5.23"
$ws.Range("E25").Value = "  +9.35%  "
$ws.Range("D26").Value = "3.306.27"
$ws.Range("E26").Value = "  -1.46%  "
$ws.Range("D27").Value = "'76.50"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").Value = "'10.84"
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").Value = "'0.0000122"
$ws.Range("E30").Value = "  +3.43%  "
$ws.Range("D31").Value = "'587.70"
$ws.Range("E31").Value = "  +12.19%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "'8.95"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").Value = "'1.51"
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("E35").Value = "  +9.50%  "
$ws.Range("D36").Value = "'1.99"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").Value = "'0.138"
$ws.Range("E37").Value = "  +15.54%  "
$ws.Range("D38").Value = "'22.72"
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("E40").Value = "  +11.82%  "
$ws.Range("D41").Value = "'0.407"
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("E42").Value = "  +14.23%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'3.07"
$ws.Range("E43").Value = "  +20.60%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "'20.76"
$ws.Range("E44").Value = "  +3.62%  "
$ws.Range("D45").Value = "'159.92"
$ws.Range("E45").Value = "  -3.21%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'187.72"
$ws.Range("E47").Value = "  -2.82%  "
$ws.Range("D48").Value = "'45.19"
$ws.Range("E48").Value = "  +5.17%  "
$ws.Range("D49").Value = "'1.34"
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("D50").Value = "'0.774"
$ws.Range("E50").Value = "  -4.76%  "
$ws.Range("D51").Value = "'26.11"
$ws.Range("E51").Value = "  +1.83%  "
